# Apply the cryptos list update (prices / 1h volume % / some re-ranked rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '27.020.77'
$ws.Range("D3").Value = '1.823.04'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.007'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4653'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.85%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3660'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07227'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.89%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8593'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.84'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07717'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.09%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.876.10'
$ws.Range("E13").Value = '  -0.05%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.325'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.69'
$ws.Range("D15").Style = "Normal"
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.490'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.25%  '
$ws.Range("E17").Value = '  -0.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008656'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.006'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.44%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.46'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.30%  '
$ws.Range("B21").Value = 'WrappedBTC'
$ws.Range("C21").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D21").Value = '26.764.21'
$ws.Range("E21").Value = '  -2.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.149'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.59%  '
$ws.Range("E23").Value = '  -1.35%  '
$ws.Range("D24").Value = '2.158.21'
$ws.Range("E24").Value = '  +3.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.65'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.839'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.10%  '
$ws.Range("E27").Value = '  -2.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.055'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.100'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.45'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08831'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.76%  '
$ws.Range("E32").Value = '  +0.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.423'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.70%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7203'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.84%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.129'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.074'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.42%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05241'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.415'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01922'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.927'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.129'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5158'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1626'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8587'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -15.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.173'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4787'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.96%  '
$ws.Range("E47").Value = '  -0.35%  '
$ws.Range("E48").Value = '  -3.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '102.55'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.42%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06241'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.616'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.32%  '
